$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Fix typo "اسم القسم" -> "أسم القسم" (shared by the Sheet1 header and Sheet3 data cell)
$ws1.Range("G1").Value = "أسم القسم"
$ws3.Range("A2").Value = "أسم القسم"

# Update selection on Sheet1 (no longer the active tab)
$ws1.Range("G1").Select()

# Update selection on Sheet3 and make it the active tab
$ws3.Range("A2").Select()
